# Update scripts with new TPM-derived NATMI ligand-receptor values for
# Rspo1-Lrp6 (OldD2 dataset): the sending/target cluster set changed
# (ECs / Inflammatory-Mac / FAPs / MuSCs / Neutrophils / Resolving-Mac)
# and every numeric metric in columns E:T was recomputed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-13 (columns A-T), in row order.
$data = @(
    @("ECs", "Rspo1", "Lrp6", "ECs", 1, 0.5, 0.724909, 1.449818, 0.5784209917375155, 0.477722171991027, 2, 1, 16.558025, 33.11605, 0.2047872600336892, 0.155585835247525, 12.003061344725, 48.0122453789, 0.118453250043895, 0.07432680314548572),
    @("ECs", "Rspo1", "Lrp6", "FAPs", 1, 0.5, 0.724909, 1.449818, 0.5784209917375155, 0.477722171991027, 3, 1, 20.31779433333333, 60.953383, 0.2512875437409564, 0.2863712008291233, 14.72855197238233, 88.37131183429399, 0.1453499902619283, 0.1368058720557674),
    @("ECs", "Rspo1", "Lrp6", "Inflammatory-Mac", 1, 0.5, 0.724909, 1.449818, 0.5784209917375155, 0.477722171991027, 3, 1, 10.405091, 31.215273, 0.1286886616182304, 0.1466556042216543, 7.542744111718999, 45.256464670314, 0.07443622327859035, 0.0700606337834251),
    @("ECs", "Rspo1", "Lrp6", "MuSCs", 1, 0.5, 0.724909, 1.449818, 0.5784209917375155, 0.477722171991027, 2, 1, 13.1587975, 26.317595, 0.1627461056112162, 0.1236453320906656, 9.5389307369275, 38.15572294771, 0.09413576380905814, 0.0590681166029046),
    @("ECs", "Rspo1", "Lrp6", "Neutrophils", 1, 0.5, 0.724909, 1.449818, 0.5784209917375155, 0.477722171991027, 3, 1, 7.303315666666667, 21.909947, 0.09032635260170116, 0.1029373190408882, 5.294239256607667, 31.76543553964601, 0.0522466584519085, 0.04917543963114639),
    @("ECs", "Rspo1", "Lrp6", "Resolving-Mac", 1, 0.5, 0.724909, 1.449818, 0.5784209917375155, 0.477722171991027, 3, 1, 13.11173766666667, 39.335213, 0.1621640763942067, 0.1848047085701436, 9.504816640205668, 57.02889984123401, 0.09379910589213526, 0.08828530677229775),
    @("Inflammatory-Mac", "Rspo1", "Lrp6", "ECs", 2, 0.6666666666666666, 0.528346, 1.585038, 0.4215790082624845, 0.522277828008973, 2, 1, 16.558025, 33.11605, 0.2047872600336892, 0.155585835247525, 8.74836627665, 52.4901976599, 0.08633400998979421, 0.08125903210203926),
    @("Inflammatory-Mac", "Rspo1", "Lrp6", "FAPs", 2, 0.6666666666666666, 0.528346, 1.585038, 0.4215790082624845, 0.522277828008973, 3, 1, 20.31779433333333, 60.953383, 0.2512875437409564, 0.2863712008291233, 10.73482536483933, 96.61342828355399, 0.1059375534790281, 0.149565328773356),
    @("Inflammatory-Mac", "Rspo1", "Lrp6", "Inflammatory-Mac", 2, 0.6666666666666666, 0.528346, 1.585038, 0.4215790082624845, 0.522277828008973, 3, 1, 10.405091, 31.215273, 0.1286886616182304, 0.1466556042216543, 5.497488209485999, 49.47739388537399, 0.05425243833964, 0.07659497043822917),
    @("Inflammatory-Mac", "Rspo1", "Lrp6", "MuSCs", 2, 0.6666666666666666, 0.528346, 1.585038, 0.4215790082624845, 0.522277828008973, 2, 1, 13.1587975, 26.317595, 0.1627461056112162, 0.1236453320906656, 6.952398023934999, 41.71438814360999, 0.0686103418021581, 0.064577215487761),
    @("Inflammatory-Mac", "Rspo1", "Lrp6", "Neutrophils", 2, 0.6666666666666666, 0.528346, 1.585038, 0.4215790082624845, 0.522277828008973, 3, 1, 7.303315666666667, 21.909947, 0.09032635260170116, 0.1029373190408882, 3.858677619220667, 34.72809857298601, 0.03807969414979265, 0.05376187940974177),
    @("Inflammatory-Mac", "Rspo1", "Lrp6", "Resolving-Mac", 2, 0.6666666666666666, 0.528346, 1.585038, 0.4215790082624845, 0.522277828008973, 3, 1, 13.11173766666667, 39.335213, 0.1621640763942067, 0.1848047085701436, 6.927534149232668, 62.347807343094, 0.06836497050207142, 0.09651940179784584)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $excelRow = $i + 2
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($excelRow, $j + 1).Value = $row[$j]
    }
}

